$d = $word.ActiveDocument

# Remove the significance-stars suffix from two estimate cells.
$d.Content.Find.Execute("42.850***", $true, $false, $false, $false, $false, $false, 1, $false, "42.850", 2)
$d.Content.Find.Execute("2.235***", $true, $false, $false, $false, $false, $false, 1, $false, "2.235", 2)

# Drop the trailing footnote row ("+ p < 0.1, * p < 0.05, ...") from the
# regression table entirely.
$t = $d.Tables.Item(1)
$t.Rows.Item($t.Rows.Count).Delete()
